# Applies the cryptos-list refresh: updates Coin/Link/Price/Volume(1h)
# cells that changed between the previous and current scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new cell text. Using NumberFormat "@" (Text) before
# assigning the value keeps numeric-looking strings (e.g. "0.200", "1.00")
# stored verbatim instead of being normalized into floating point numbers,
# then resetting the style back to "Normal" afterwards so the cell keeps the
# workbook default formatting (only the textual content changes).
$updates = [ordered]@{
    "D2" = "46.555.76"
    "E2" = "  +5.49%  "
    "D3" = "2.308.90"
    "E3" = "  +5.03%  "
    "E4" = "  -0.87%  "
    "D5" = "299.44"
    "E5" = "  +1.03%  "
    "D6" = "96.97"
    "E6" = "  +8.27%  "
    "E7" = "  -0.18%  "
    "D8" = "0.999"
    "E8" = "  -0.66%  "
    "E9" = "  +9.04%  "
    "D10" = "35.35"
    "E10" = "  +7.08%  "
    "D11" = "0.0797"
    "E11" = "  +2.47%  "
    "D12" = "7.34"
    "E12" = "  +8.48%  "
    "E13" = "  +1.17%  "
    "D14" = "2.659.10"
    "E14" = "  +4.57%  "
    "D15" = "2.308.10"
    "E15" = "  +4.11%  "
    "D16" = "13.97"
    "E16" = "  +6.10%  "
    "D17" = "0.822"
    "E17" = "  +6.62%  "
    "D18" = "46.522.11"
    "E18" = "  +6.00%  "
    "D19" = "13.10"
    "E19" = "  +20.24%  "
    "D20" = "0.0₃0937"
    "E20" = "  +5.50%  "
    "D21" = "6.13"
    "E21" = "  +5.10%  "
    "D22" = "66.86"
    "E22" = "  +5.61%  "
    "D23" = "248.39"
    "E23" = "  +8.31%  "
    "D24" = "2.91"
    "E24" = "  +5.28%  "
    "E25" = "  +8.59%  "
    "E26" = "  -0.33%  "
    "D27" = "42.40"
    "E27" = "  +18.58%  "
    "D28" = "2.23"
    "E28" = "  +0.90%  "
    "D29" = "9.82"
    "E29" = "  +6.29%  "
    "D30" = "20.09"
    "E30" = "  +5.68%  "
    "D31" = "5.77"
    "E31" = "  +8.58%  "
    "B32" = "Monero"
    "C32" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D32" = "147.03"
    "E32" = "  +0.75%  "
    "B33" = "Hedera"
    "C33" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D33" = "0.0802"
    "E33" = "  +8.83%  "
    "D34" = "2.61"
    "E34" = "  +4.30%  "
    "D35" = "3.12"
    "E35" = "  +8.67%  "
    "D36" = "0.112"
    "E36" = "  +9.14%  "
    "E37" = "  +1.78%  "
    "E38" = "  +9.18%  "
    "B39" = "Celestia"
    "C39" = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
    "D39" = "15.14"
    "E39" = "  +11.87%  "
    "B40" = "RenderToken"
    "C40" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D40" = "3.99"
    "E40" = "  +13.89%  "
    "D41" = "3.44"
    "E41" = "  +10.80%  "
    "E42" = "  +8.07%  "
    "D43" = "0.999"
    "E43" = "  -0.82%  "
    "B44" = "Stacks"
    "C44" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D44" = "1.97"
    "E44" = "  +19.80%  "
    "B45" = "Maker"
    "C45" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D45" = "1.836.26"
    "E45" = "  +4.49%  "
    "D46" = "90.06"
    "E46" = "  +21.42%  "
    "D47" = "0.200"
    "E47" = "  +15.95%  "
    "D48" = "72.67"
    "E48" = "  +4.71%  "
    "D49" = "4.96"
    "E49" = "  +12.71%  "
    "D50" = "97.43"
    "E50" = "  +5.90%  "
    "D51" = "54.21"
    "E51" = "  +9.22%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
